$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D4,D5,D6,D10,D11,D12,D14,D17,D19,D20,D21,D22,D24,D25,D26,D28,D29,D30,D31,D32,D33,D34,D35,D36,D38,D39,D40,D41,D42,D43,D44,D45,D46,D49,D50,D51').NumberFormat = '@'

$ws.Range('D2').Value = '61.736.78'
$ws.Range('E2').Value = '  +0.62%  '

$ws.Range('D3').Value = '3.454.61'
$ws.Range('E3').Value = '  +0.65%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').Value = '578.87'
$ws.Range('E5').Value = '  +0.71%  '

$ws.Range('D6').Value = '146.15'
$ws.Range('E6').Value = '  +4.49%  '

$ws.Range('D7').Value = '3.454.79'
$ws.Range('E7').Value = '  +0.66%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('E9').Value = '  +1.69%  '

$ws.Range('D10').Value = '7.68'
$ws.Range('E10').Value = '  -0.94%  '

$ws.Range('D11').Value = '0.128'
$ws.Range('E11').Value = '  +3.81%  '

$ws.Range('D12').Value = '0.392'
$ws.Range('E12').Value = '  +2.56%  '

$ws.Range('D13').Value = '4.043.42'
$ws.Range('E13').Value = '  +0.61%  '

$ws.Range('D14').Value = '28.79'

$ws.Range('E15').Value = '  -0.47%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.495.72'
$ws.Range('E16').Value = '  +1.97%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.0000175'
$ws.Range('E17').Value = '  +1.20%  '

$ws.Range('D18').Value = '61.789.01'
$ws.Range('E18').Value = '  +0.52%  '

$ws.Range('D19').Value = '6.39'

$ws.Range('D20').Value = '14.37'
$ws.Range('E20').Value = '  +2.94%  '

$ws.Range('D21').Value = '9.45'
$ws.Range('E21').Value = '  +0.09%  '

$ws.Range('D22').Value = '406.26'
$ws.Range('E22').Value = '  +5.71%  '

$ws.Range('E23').Value = '  +2.17%  '

$ws.Range('D24').Value = '74.36'
$ws.Range('E24').Value = '  +3.63%  '

$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.80%  '

$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '0.0000124'
$ws.Range('E26').Value = '  +1.06%  '

$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '3.592.34'
$ws.Range('E27').Value = '  +1.36%  '

$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').Value = '0.183'
$ws.Range('E28').Value = '  +3.75%  '

$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '7.65'
$ws.Range('E29').Value = '  +1.54%  '

$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.21%  '

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '8.27'
$ws.Range('E31').Value = '  +1.61%  '

$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '1.48'
$ws.Range('E32').Value = '  -7.03%  '

$ws.Range('D33').Value = '2.20'
$ws.Range('E33').Value = '  +1.96%  '

$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  -0.05%  '

$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').Value = '24.05'
$ws.Range('E35').Value = '  +1.20%  '

$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').Value = '7.09'
$ws.Range('E36').Value = '  +2.09%  '

$ws.Range('B37').Value = 'RenzoRestakedETH'
$ws.Range('C37').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D37').Value = '3.478.64'
$ws.Range('E37').Value = '  +0.78%  '

$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '1.57'
$ws.Range('E38').Value = '  +0.44%  '

$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = '5.17'
$ws.Range('E39').Value = '  +0.21%  '

$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').Value = '166.93'
$ws.Range('E40').Value = '  +0.62%  '

$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '0.0797'
$ws.Range('E41').Value = '  +2.45%  '

$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '27.38'
$ws.Range('E42').Value = '  +3.75%  '

$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '0.805'
$ws.Range('E43').Value = '  +2.59%  '

$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = '1.75'
$ws.Range('E44').Value = '  +0.35%  '

$ws.Range('D45').Value = '4.54'
$ws.Range('E45').Value = '  +2.62%  '

$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '42.48'
$ws.Range('E46').Value = '  +0.68%  '

$ws.Range('E47').Value = '  -0.07%  '

$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.610.75'
$ws.Range('E48').Value = '  +1.21%  '

$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').Value = '1.16'
$ws.Range('E49').Value = '  -1.94%  '

$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '7.00'
$ws.Range('E50').Value = '  +2.73%  '

$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '23.25'
$ws.Range('E51').Value = '  -2.67%  '

$ws.Range('D4,D5,D6,D10,D11,D12,D14,D17,D19,D20,D21,D22,D24,D25,D26,D28,D29,D30,D31,D32,D33,D34,D35,D36,D38,D39,D40,D41,D42,D43,D44,D45,D46,D49,D50,D51').Style = 'Normal'
